$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = @(-0.006450402548902858,-0.019333790650941686,-0.0119863214902941,-0.016207721772448242,-0.005859896842051722,-0.036597770200784,-0.01169790902106695,-0.014110421963665932,-0.010216880691312283,-0.010193501179635595,-0.013228902677560784,-0.02069906381626862,-0.008678211792776988,-0.008775447458140921,-0.015883962289262508,-0.021950572130098878,-0.012828617042535273,-0.028631710755552448,-0.01854574186869594,-0.033333838903449854,-0.010495163767167227,-0.0087718176530525,-0.012000092562800069,-0.012919420119327365,-0.012577950587654114,-0.00620344735404875,-0.0206675444796317,-0.015964217160853788,-0.011111674450826566,-0.010420416214200824,-0.015299933778063893,-0.029388368237555612,-0.01198312298118725,-0.00864524735068022,-0.01053650889280548,-0.006444881452741182,-0.008221331269418108,-0.011506005474051242,-0.01780633661093218,-0.026768957523729466,-0.011488410332588133,-0.005211564097503908,-0.004617839530774104,-0.015676637159600743,-0.021072937624592888,-0.007706164770793434,-0.022808428278693914,-0.0058129226812853926,-0.020065522596864293,-0.014665906139896096,-0.014972742081440086,-0.01890375059733526,-0.020650695321250266,-0.03141999179928992,-0.009660841088507278,-0.014923627255984923,-0.015835564572503225,-0.0287608265973477,-0.01035577165001391,-0.01191482878189001,-0.006394800896915119,-0.004801876170490373,-0.009768649509042951,-0.009307551278584203,-0.00782172668745368,-0.007972971160012143,-0.013040982048442944,-0.018312163998991285,-0.024533619239460704,-0.030561830880395137,-0.02009027984570543,-0.009130253078267765,-0.021343920016180434,-0.010832143285207698,-0.037009788860856624,-0.01893267733397676,-0.007399550693435029,-0.01114469259357677,-0.014198038996972693,-0.02246470641695452,-0.009161854648721171,-0.013818485208436294,-0.005920727151504787,-0.006938676902238857,-0.01849210047179895,-0.014125327532368318,-0.019590092891750864,-0.01750813193372548,-0.031158406788846425,-0.0112551168066724,-0.013307974771696252,-0.015499828232168004,-0.016053430893945873,-0.013532107186406739,-0.03448837330249046,-0.012261999536137203,-0.020832137760713586,-0.014416657991352047,-0.01907986811211293,-0.007441710427783249)
$row2 = @(-0.00165217313143688,-0.004887028213155732,-0.002920789039641706,-0.004017954131925455,-0.0014033230078909734,-0.00918523752832961,-0.0028610648200679202,-0.0035208111354968438,-0.002557554258571318,-0.002618873542941332,-0.003268914715008224,-0.005178789766817468,-0.002208062901277619,-0.0022189817114851907,-0.00401519515537542,-0.005519992912638728,-0.0032378934701796642,-0.007210714115826268,-0.0045881732976117895,-0.00837030259830774,-0.0026476250718255917,-0.0021071514086722777,-0.0029493820864654103,-0.0033032505499956994,-0.003185630021106592,-0.0015149163330558828,-0.0052478276312444005,-0.0039355338786642696,-0.0027694454629834725,-0.002624331431743447,-0.0037780447277469273,-0.007428364845114439,-0.0029772350701304058,-0.0021993827629737833,-0.002622621468732593,-0.0015746497883529554,-0.0020047454893189352,-0.002861638343670881,-0.004389358906293836,-0.00675158031299243,-0.002836903517528909,-0.0012848017510258215,-0.001100604624470266,-0.003958912219331106,-0.005248041873585464,-0.0019147354599405068,-0.005720546729175373,-0.0014613393204559756,-0.0049668252620808356,-0.0036395847089472247,-0.0037425115050218794,-0.004724170281506176,-0.005218520171516493,-0.007946676218920798,-0.002454583142046469,-0.0036950289180331156,-0.003976591901794242,-0.007230135291535117,-0.0025392163681093056,-0.0030033973956990145,-0.0016222055843765119,-0.00122864401393714,-0.0024411507800483003,-0.0022838506727168146,-0.0019659006979968695,-0.001997541498459083,-0.0032375241230661973,-0.004588288301524842,-0.006195820622143616,-0.007736341512280435,-0.004963450399777967,-0.002259242393405192,-0.00538198045135544,-0.0027073030797323443,-0.009293119234790547,-0.004862328461296222,-0.0018462693300767392,-0.0028089857094220558,-0.0035007997558170975,-0.005584669706900814,-0.0023278155563475514,-0.00341305094438991,-0.0014593263032712303,-0.001700708860242543,-0.004726680282313076,-0.0035888113516663266,-0.0050321337425578054,-0.004390788441401899,-0.007790174372951692,-0.0027703743336284226,-0.0034071835867211485,-0.003979717204567019,-0.003989911973944423,-0.003292458477051593,-0.008665059748638715,-0.0030212371232287865,-0.005316146895827898,-0.0035824140823284244,-0.004741249721471677,-0.0018115421036353292)

$arr1 = New-Object 'object[,]' 1,100
$arr2 = New-Object 'object[,]' 1,100

for ($i = 0; $i -lt 100; $i++) {
    $arr1[0, $i] = $row1[$i]
    $arr2[0, $i] = $row2[$i]
}

$ws.Range("A1:CV1").Value2 = $arr1
$ws.Range("A2:CV2").Value2 = $arr2
